# Remove the "Emails / Investigate having the from address set to the
# user's email address" backlog item (row 1). Deleting the entire row
# shifts every row below it up by one and drops the two now-unused
# shared strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

# Match the author's saved cursor position (B6) after the edit.
$ws.Range("B6").Select()
